$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet updates ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Beckley Pocahontas Coal Mine, United States, M0990, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$wsAbout.Range("A6").Value = $newCitation

# --- "Boundaries and methane sources" sheet updates ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 7; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19 (build_version)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
